$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "xNome"
$ws.Range("H2").Value = "CAFE RANCHEIRO AGRO INDUSTRIAL LTDA"
